# Generate Report for Handoff
# Updates the handoff timestamps and priority markers for the rows that
# were just generated/handed off (rows 7, 9, 10, 12, 13, 14 on each sheet).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column G
    $overview.Range("G$r").Value = "2016-09-05 20:25:53"

    # zh-cn sheet: "Latest Handoff Datetime" column H, "Priority" column E
    $zhcn.Range("H$r").Value = "2016-09-05 20:25:47"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: "Latest Handoff Datetime" column H, "Priority" column E
    $dede.Range("H$r").Value = "2016-09-05 20:25:53"
    $dede.Range("E$r").Value = "ht"
}
